$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include from Observation code")

$ws.Range("B2").Value = "Serum Creatine Kinase"
$ws.Range("B3").Value = "Muscle Biopsy"
$ws.Range("B7").Value = "GCN Repeat Testing (Oculopharyngeal Muscular Dystrophy)"
$ws.Range("B8").Value = "Deletions and Duplications Testing (Duchenne and Becker Dystrophies)"
